$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 57, shifting existing rows 57:90 down to 58:91
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new record's data
$ws.Range("A57").Value = 3
$ws.Range("B57").Value = "Femacal de La Calera"
$ws.Range("C57").Value = "Coquimbo"
$ws.Range("D57").Value = 44510
$ws.Range("E57").Value = 5
$ws.Range("F57").Value = 100112026
$ws.Range("G57").Value = "Haba"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 90
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 7500
$ws.Range("M57").Value = 7278
$ws.Range("N57").Value = "$/malla 25 kilos"
$ws.Range("O57").Value = "Provincia de Quillota"
$ws.Range("P57").Value = 291
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"

# Apply the same date number format as the other "Fecha" cells in column D
$ws.Range("D57").NumberFormat = "YYYY-MM-DD HH:MM:SS"
